# Pixel Buffer Memory Calculations.xlsx -- add a "Total Memory" worksheet
# (Sheet2) that rolls up the model (points) memory size together with the
# existing pixel/z-buffer image size, and make it the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$fmtInt = '_(* #,##0_);_(* \(#,##0\);_(* "-"??_);_(@_)'
$fmt5dp = '_(* #,##0.00000_);_(* \(#,##0.00000\);_(* "-"??_);_(@_)'

# ---------------------------------------------------------------------
# Row 1-2: Width / Height
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "Width:"
$ws2.Range("A1").Font.Italic = $true
$ws2.Range("A1").HorizontalAlignment = -4152  # xlRight
$ws2.Range("B1").Value = 1500
$ws2.Range("B1").NumberFormat = $fmtInt

$ws2.Range("A2").Value = "Height:"
$ws2.Range("A2").Font.Italic = $true
$ws2.Range("A2").HorizontalAlignment = -4152
$ws2.Range("B2").Value = 1500
$ws2.Range("B2").NumberFormat = $fmtInt

# ---------------------------------------------------------------------
# Row 3: Total Points
# ---------------------------------------------------------------------
$ws2.Range("A3").Value = "Total Points:"
$ws2.Range("A3").Font.Underline = $true
$ws2.Range("A3").HorizontalAlignment = -4152
$ws2.Range("B3").Formula = "=B1*B2"
$ws2.Range("B3").Font.Underline = $true
$ws2.Range("B3").NumberFormat = $fmtInt

# Row 4: blank spacer (still carries the underline style of row 3)
$ws2.Range("A4").Font.Underline = $true
$ws2.Range("A4").HorizontalAlignment = -4152
$ws2.Range("B4").Font.Underline = $true
$ws2.Range("B4").NumberFormat = $fmtInt

# ---------------------------------------------------------------------
# Row 5: Subpixel Grid
# ---------------------------------------------------------------------
$ws2.Range("A5").Value = "Subpixel Grid:"
$ws2.Range("A5").Font.Italic = $true
$ws2.Range("A5").HorizontalAlignment = -4152
$ws2.Range("B5").Value = 4
$ws2.Range("B5").NumberFormat = $fmtInt

# ---------------------------------------------------------------------
# Row 6: Total Samples (model byte size) - uses the new underline font
# ---------------------------------------------------------------------
$ws2.Range("A6").Value = "Total Samples:"
$ws2.Range("A6").Font.Underline = $true
$ws2.Range("A6").HorizontalAlignment = -4152
$ws2.Range("B6").Formula = "=B3*B5"
$ws2.Range("B6").Font.Underline = $true
$ws2.Range("B6").NumberFormat = $fmtInt
$ws2.Rows.Item(6).RowHeight = 17.25

# ---------------------------------------------------------------------
# Rows 7-13: point-structure byte breakdown table
# ---------------------------------------------------------------------
$ws2.Range("A7").HorizontalAlignment = -4152

$ws2.Range("A8").HorizontalAlignment = -4152
$ws2.Range("B8").Value = "Type"
$ws2.Range("B8").HorizontalAlignment = -4108  # xlCenter
$ws2.Range("C8").Value = "Bytes"
$ws2.Range("C8").HorizontalAlignment = -4108

$ws2.Range("A9").Value = "Elevation:"
$ws2.Range("A9").HorizontalAlignment = -4152
$ws2.Range("B9").Value = "double"
$ws2.Range("B9").Font.Italic = $true
$ws2.Range("C9").Value = 8
$ws2.Range("C9").Font.Name = "Calibri"

$ws2.Range("A10").Value = "Normal:"
$ws2.Range("A10").HorizontalAlignment = -4152
$ws2.Range("B10").Value = "double[3]"
$ws2.Range("B10").Font.Italic = $true
$ws2.Range("C10").Value = 24
$ws2.Range("C10").Font.Name = "Calibri"

$ws2.Range("A11").Value = "Dot Product:"
$ws2.Range("A11").HorizontalAlignment = -4152
$ws2.Range("B11").Value = "double"
$ws2.Range("B11").Font.Italic = $true
$ws2.Range("C11").Value = 8
$ws2.Range("C11").Font.Name = "Calibri"

$ws2.Range("A12").Value = "RGBA:"
$ws2.Range("A12").HorizontalAlignment = -4152
$ws2.Range("B12").Value = "int"
$ws2.Range("B12").Font.Italic = $true
$ws2.Range("C12").Value = 4
$ws2.Range("C12").Font.Name = "Calibri"

$ws2.Range("A13").Value = "Size:"
$ws2.Range("A13").Font.Underline = $true
$ws2.Range("A13").HorizontalAlignment = -4152
$ws2.Range("C13").Formula = "=SUM(C9:C12)"
$ws2.Range("C13").Font.Underline = $true

$ws2.Range("A14").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# Rows 15-16: Total Model Size
# ---------------------------------------------------------------------
$ws2.Range("A15").HorizontalAlignment = -4152
$ws2.Range("B15").Value = "Bytes:"
$ws2.Range("B15").Font.Bold = $true
$ws2.Range("C15").Value = "Megabytes:"
$ws2.Range("C15").Font.Bold = $true
$ws2.Range("D15").Value = "Gigabytes:"
$ws2.Range("D15").Font.Bold = $true

$ws2.Range("A16").Value = "Total Model Size:"
$ws2.Range("A16").Font.Bold = $true
$ws2.Range("A16").Font.Underline = $true
$ws2.Range("A16").HorizontalAlignment = -4152
$ws2.Range("B16").Formula = "=B6*C13"
$ws2.Range("B16").Font.Bold = $true
$ws2.Range("B16").Font.Underline = $true
$ws2.Range("B16").NumberFormat = $fmtInt
$ws2.Range("C16").Formula = "=B16/POWER(2, 20)"
$ws2.Range("C16").Font.Bold = $true
$ws2.Range("C16").Font.Underline = $true
$ws2.Range("C16").NumberFormat = $fmt5dp
$ws2.Range("D16").Formula = "=B16/POWER(2, 30)"
$ws2.Range("D16").Font.Bold = $true
$ws2.Range("D16").Font.Underline = $true
$ws2.Range("D16").NumberFormat = $fmt5dp

# ---------------------------------------------------------------------
# Rows 19-20: Bytes per color / z-buffer sample
# ---------------------------------------------------------------------
$ws2.Range("A19").Value = "Bytes Per Color Sample:"
$ws2.Range("A19").Font.Underline = $true
$ws2.Range("A19").HorizontalAlignment = -4152
$ws2.Range("B19").Value = 4
$ws2.Range("B19").Font.Underline = $true
$ws2.Range("B19").NumberFormat = $fmtInt

$ws2.Range("A20").Value = "Bytes Per Z-Buffer Sample:"
$ws2.Range("A20").Font.Underline = $true
$ws2.Range("A20").HorizontalAlignment = -4152
$ws2.Range("B20").Value = 8
$ws2.Range("B20").Font.Underline = $true
$ws2.Range("B20").NumberFormat = $fmtInt

$ws2.Range("B21").NumberFormat = $fmtInt

# ---------------------------------------------------------------------
# Rows 22-26: Total Image Size (color + z buffers)
# ---------------------------------------------------------------------
$ws2.Range("B22").Value = "Bytes:"
$ws2.Range("B22").Font.Bold = $true
$ws2.Range("B22").NumberFormat = $fmtInt
$ws2.Range("C22").Value = "Megabytes:"
$ws2.Range("C22").Font.Bold = $true
$ws2.Range("D22").Value = "Gigabytes:"
$ws2.Range("D22").Font.Bold = $true

$ws2.Range("A23").Value = "Total Color Buffer:"
$ws2.Range("A23").Font.Underline = $true
$ws2.Range("A23").HorizontalAlignment = -4152
$ws2.Range("B23").Formula = "=B6 * B19"
$ws2.Range("B23").Font.Underline = $true
$ws2.Range("B23").NumberFormat = $fmtInt
$ws2.Range("C23").Formula = "=B23/POWER(2, 20)"
$ws2.Range("C23").Font.Underline = $true
$ws2.Range("C23").NumberFormat = $fmt5dp
$ws2.Range("D23").Formula = "=B23/POWER(2, 30)"
$ws2.Range("D23").Font.Underline = $true
$ws2.Range("D23").NumberFormat = $fmt5dp

$ws2.Range("A24").Value = "Total Z-Buffer:"
$ws2.Range("A24").Font.Underline = $true
$ws2.Range("A24").HorizontalAlignment = -4152
$ws2.Range("B24").Formula = "=B6*B20"
$ws2.Range("B24").Font.Underline = $true
$ws2.Range("B24").NumberFormat = $fmtInt
$ws2.Range("C24").Formula = "=B24/POWER(2, 20)"
$ws2.Range("C24").Font.Underline = $true
$ws2.Range("C24").NumberFormat = $fmt5dp
$ws2.Range("D24").Formula = "=B24/POWER(2, 30)"
$ws2.Range("D24").Font.Underline = $true
$ws2.Range("D24").NumberFormat = $fmt5dp

$ws2.Range("B25").NumberFormat = $fmtInt
$ws2.Range("D25").NumberFormat = $fmt5dp

$ws2.Range("A26").Value = "Total Image Size:"
$ws2.Range("A26").Font.Bold = $true
$ws2.Range("A26").Font.Underline = $true
$ws2.Range("A26").HorizontalAlignment = -4152
$ws2.Range("B26").Formula = "=B23+B24"
$ws2.Range("B26").Font.Bold = $true
$ws2.Range("B26").Font.Underline = $true
$ws2.Range("B26").NumberFormat = $fmtInt
$ws2.Range("C26").Formula = "=B26/POWER(2, 20)"
$ws2.Range("C26").Font.Bold = $true
$ws2.Range("C26").Font.Underline = $true
$ws2.Range("C26").NumberFormat = $fmt5dp
$ws2.Range("D26").Formula = "=B26/POWER(2, 30)"
$ws2.Range("D26").Font.Bold = $true
$ws2.Range("D26").Font.Underline = $true
$ws2.Range("D26").NumberFormat = $fmt5dp

# ---------------------------------------------------------------------
# Rows 29-30: Total Memory Size (model + image)
# ---------------------------------------------------------------------
$ws2.Range("B29").Value = "Bytes:"
$ws2.Range("B29").Font.Bold = $true
$ws2.Range("C29").Value = "Megabytes:"
$ws2.Range("C29").Font.Bold = $true
$ws2.Range("D29").Value = "Gigabytes:"
$ws2.Range("D29").Font.Bold = $true

$ws2.Range("A30").Value = "Total Memory Size:"
$ws2.Range("A30").Font.Bold = $true
$ws2.Range("A30").Font.Underline = $true
$ws2.Range("A30").HorizontalAlignment = -4152
$ws2.Range("B30").Formula = "=B16+B26"
$ws2.Range("B30").Font.Bold = $true
$ws2.Range("B30").Font.Underline = $true
$ws2.Range("B30").NumberFormat = $fmtInt
$ws2.Range("C30").Formula = "=B30/POWER(2, 20)"
$ws2.Range("C30").Font.Bold = $true
$ws2.Range("C30").Font.Underline = $true
$ws2.Range("C30").NumberFormat = $fmt5dp
$ws2.Range("D30").Formula = "=B30/POWER(2, 30)"
$ws2.Range("D30").Font.Bold = $true
$ws2.Range("D30").Font.Underline = $true
$ws2.Range("D30").NumberFormat = $fmt5dp

# ---------------------------------------------------------------------
# Column widths (auto-fit approximations)
# ---------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 24.14
$ws2.Columns.Item(2).ColumnWidth = 13.42
$ws2.Columns.Item(3).ColumnWidth = 11.6
$ws2.Columns.Item(4).ColumnWidth = 10.2

# ---------------------------------------------------------------------
# Views / selection: Sheet2 becomes the active tab
# ---------------------------------------------------------------------
$ws1.Range("A10:D17").Select()

$ws2.Activate()
$ws2.Application.ActiveWindow.Zoom = 85
$ws2.Range("C18").Select()
